$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 (Mogastone Level 4): price bump + "Includes Sink" No -> Yes ---
$ws.Range("D12").Value = 3358.62
$ws.Range("F12").Value = "Yes"

# --- Row 13 (Mogastone Level 5): price bump (now shown with 2-decimal format) + sink Yes ---
$ws.Range("D13").Value = 3843.97
$ws.Range("D13").NumberFormat = "0.00"
$ws.Range("F13").Value = "Yes"

# --- Row 14 (Mogastone Level 6): price bump + sink Yes ---
$ws.Range("D14").Value = 4188.6000000000004
$ws.Range("F14").Value = "Yes"

# --- Column D was widened to fit the new decimal pricing ---
$ws.Columns("D").ColumnWidth = 11.3333333333333

# --- Leave the selection where the last edit was made ---
$ws.Range("F14").Select()
